$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 contents (set B2 first so the shared-string table ends up
# ordered the same way as in the target workbook):
# B2: "That learner has knowledge of the fgr loop. You should ask the learner complex questions"
#     -> "That learner has knowledge of the for loop. You should ask the learner complex questions" (typo fix)
# A2: "A learner can write the for loop and knows the syntax of the for loop" -> "I can write a for loop"
$ws.Range("B2").Value = "That learner has knowledge of the for loop. You should ask the learner complex questions"
$ws.Range("A2").Value = "I can write a for loop"

# Give column K (11) the same explicit custom width as columns D:J
$ws.Columns("K:K").ColumnWidth = $ws.Columns("D:D").ColumnWidth

# Move the active selection to B6
$ws.Range("B6").Select()
